# Scen_NODESTEC.xlsx cleanup edit
# Adds a new "START" / 2100 / "EUWINOFV01" row to the INS sheet (mirrors the
# existing START/2100/IMPELC-DESTEC row) and leaves the active selection on
# the cell just below the new entry, matching the author's manual edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INS")

$ws.Range("D5").Value = "START"
$ws.Range("H5").Value = 2100
$ws.Range("J5").Value = "EUWINOFV01"

[void]$ws.Range("J6").Select()
